# Update cd_fcv sheet: split each metric column into "mean"/"std" pairs,
# rename CART -> DTREE, drop the NB row, and refresh all the numeric data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Make room for the three new "std" columns -------------------------
# Current layout: B=Algorithm, C=State Based, D=Non State, E=One Sided
# Target layout:  B=Algorithm, C=State Based mean, D=State Based std,
#                 E=Non State mean, F=Non State std,
#                 G=One Sided mean, H=One Sided std
$ws.Columns("D:D").Insert()
$ws.Columns("F:F").Insert()
$ws.Columns("H:H").Insert()

# Copy the header style (bold / centered / bordered) from an existing
# header cell onto the freshly inserted header cells.
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("H1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- 2. Header row ----------------------------------------------------------
$ws.Range("B1").Value = "Algorithm"
$ws.Range("C1").Value = "State Based mean"
$ws.Range("D1").Value = "State Based std"
$ws.Range("E1").Value = "Non State mean"
$ws.Range("F1").Value = "Non State std"
$ws.Range("G1").Value = "One Sided mean"
$ws.Range("H1").Value = "One Sided std"

# --- 3. Data rows -------------------------------------------------------
$ws.Range("B2").Value = "LR"
$ws.Range("C2").Value = 0.8025081255048028
$ws.Range("D2").Value = 0.02356802553124964
$ws.Range("E2").Value = 0.6714808653926208
$ws.Range("F2").Value = 0.02608003172789519
$ws.Range("G2").Value = 0.7728582787472581
$ws.Range("H2").Value = 0.01657077722819842

$ws.Range("B3").Value = "LDA"
$ws.Range("C3").Value = 0.7974741785752899
$ws.Range("D3").Value = 0.01992704435182157
$ws.Range("E3").Value = 0.647180368849116
$ws.Range("F3").Value = 0.02407491536831224
$ws.Range("G3").Value = 0.7701251337001042
$ws.Range("H3").Value = 0.01740480239606836

$ws.Range("B4").Value = "KNN"
$ws.Range("C4").Value = 0.7382144119644479
$ws.Range("D4").Value = 0.0317805925855184
$ws.Range("E4").Value = 0.6382159588736129
$ws.Range("F4").Value = 0.02306334582189203
$ws.Range("G4").Value = 0.7286168171583105
$ws.Range("H4").Value = 0.01362171461298567

$ws.Range("B5").Value = "DTREE"
$ws.Range("C5").Value = 0.765062159800849
$ws.Range("D5").Value = 0.02883513785210558
$ws.Range("E5").Value = 0.6525911950291723
$ws.Range("F5").Value = 0.03895014491171871
$ws.Range("G5").Value = 0.7326965321209434
$ws.Range("H5").Value = 0.02277228366219797

$ws.Range("B6").Value = "RTREE"
$ws.Range("C6").Value = 0.7439006295713522
$ws.Range("D6").Value = 0.03012637725099197
$ws.Range("E6").Value = 0.5317685434439352
$ws.Range("F6").Value = 0.0228584746169964
$ws.Range("G6").Value = 0.7098032611002263
$ws.Range("H6").Value = 0.01657264814403177

$ws.Range("B7").Value = "XTREE"
$ws.Range("C7").Value = 0.822885070482641
$ws.Range("D7").Value = 0.03324605361348739
$ws.Range("E7").Value = 0.6967736195355544
$ws.Range("F7").Value = 0.03429055084177973
$ws.Range("G7").Value = 0.7834316660479297
$ws.Range("H7").Value = 0.01964014027279924

$ws.Range("B8").Value = "SVM"
$ws.Range("C8").Value = 0.8071183292856168
$ws.Range("D8").Value = 0.02190890165737961
$ws.Range("E8").Value = 0.7155770806149129
$ws.Range("F8").Value = 0.02325134259473407
$ws.Range("G8").Value = 0.7918632246472048
$ws.Range("H8").Value = 0.01792339290204064

# --- 4. Drop the old row 9 (NB), which no longer exists in the new table ---
$ws.Rows("9:9").Delete()
